$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage purely-numeric-looking replacement text (e.g.
# "0.17", "-0.01") through a text formula. A direct .Value assignment of
# such a string gets auto-coerced into a numeric cell by Excel's type
# inference, and forcing text via a leading apostrophe / NumberFormat="@"
# stamps a quotePrefix / text-format style onto the cell - neither of
# which the target file has (its numeric-looking entries stay plain
# <t>-type shared strings with the default, unstyled cell format). Routing
# the text through Copy / PasteSpecial(xlPasteValues) keeps the
# destination a plain shared-string cell with its original style intact.
$scratch = $ws.Range("Z100")

function Set-TextValue {
  param($Cell, $Text)

  if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
    $scratch.Formula = "=""" + $Text + """"
    $scratch.Copy()
    $ws.Range($Cell).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
    $excel.CutCopyMode = $false
  } else {
    $ws.Range($Cell).Value = $Text
  }
}

Set-TextValue "B2" "0.17"
Set-TextValue "B3" "-0.01"
Set-TextValue "B4" "-0.09"
Set-TextValue "C2" "44.29***"
Set-TextValue "C3" "2.21***"
Set-TextValue "C4" "0.98"
Set-TextValue "D2" "-0.89"
Set-TextValue "D3" "0.46***"
Set-TextValue "D4" "0.82*"
